$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D contain numeric-looking strings (e.g. "44.182.46") that must remain
# plain text, matching the original inlineStr cell type. Force text storage by
# temporarily applying a text number format, then reset the style back to Normal
# so no stray style index is left attached to the cell.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '44.182.46'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.432.07'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.96%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.58'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.513'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.53%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -0.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.47'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.62%  '
$ws.Range('E11').Value = '  +1.25%  '
$ws.Range('E12').Value = '  +2.73%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.73'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.56%  '
$ws.Range('E14').Value = '  +2.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.807.90'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.458.11'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.832'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.72%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '44.176.78'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.36'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.46'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.79%  '
$ws.Range('E21').Value = '  +1.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.71'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('E23').Value = '  +3.46%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '240.52'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.48'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.93%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.35'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.35'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.55'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.41%  '
$ws.Range('E30').Value = '  +4.10%  '
$ws.Range('E31').Value = '  +15.91%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.69'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +8.79%  '
$ws.Range('E33').Value = '  +1.50%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0760'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.18%  '
$ws.Range('E36').Value = '  +3.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.54'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '129.82'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +22.42%  '
$ws.Range('E39').Value = '  +4.53%  '
$ws.Range('E40').Value = '  -0.79%  '
$ws.Range('E41').Value = '  +0.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.25'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.30%  '
$ws.Range('E43').Value = '  +2.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.963.17'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.27%  '
$ws.Range('E45').Value = '  +2.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.89'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.39'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.18%  '
$ws.Range('E48').Value = '  +8.79%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '53.43'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.97%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.73'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.28%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.16'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.43%  '
